$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 1
    6 = 2
    7 = 1
    8 = 0
    9 = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 3
    33 = 3
    34 = 2
    35 = 1
    36 = 2
    37 = 2
    38 = 2
    39 = 2
    40 = 0
    41 = 1
    42 = 2
    43 = 0
    44 = 0
    45 = 1
    46 = 1
    47 = 1
    48 = 2
    49 = 0
    50 = 1
    51 = 0
    52 = 3
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    58 = 2
    59 = 1
    60 = 0
    61 = 0
    62 = 0
    63 = 7
    64 = 1
    65 = 2
    66 = 1
    67 = 2
    68 = 0
    69 = 2
    70 = 1
    71 = 0
    72 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

